# Auto-generated script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.283.09'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '3.077.12'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''234.64'
$ws.Range('E5').Value = '  -3.03%  '
$ws.Range('D6').Value = '''608.18'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('D8').Value = '''0.377'
$ws.Range('E8').Value = '  -5.14%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '''0.805'
$ws.Range('E10').Value = '  +6.88%  '
$ws.Range('D11').Value = '3.072.40'
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('D13').Value = '93.896.95'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '''0.0000239'
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('D15').Value = '''33.62'
$ws.Range('E15').Value = '  -2.47%  '
$ws.Range('D16').Value = '''5.31'
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').Value = '3.645.36'
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '3.049.08'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').Value = '''3.53'
$ws.Range('E19').Value = '  -6.81%  '
$ws.Range('D20').Value = '''14.38'
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').Value = '''5.67'
$ws.Range('E21').Value = '  -2.26%  '
$ws.Range('D22').Value = '''438.00'
$ws.Range('E22').Value = '  -2.14%  '
$ws.Range('D23').Value = '''8.77'
$ws.Range('E23').Value = '  -6.07%  '
$ws.Range('D24').Value = '''0.0000188'
$ws.Range('E24').Value = '  -7.75%  '
$ws.Range('D25').Value = '''8.35'
$ws.Range('E25').Value = '  +5.93%  '
$ws.Range('D26').Value = '''5.48'
$ws.Range('E26').Value = '  -5.33%  '
$ws.Range('D27').Value = '''84.49'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('D28').Value = '''11.83'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').Value = '3.225.46'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').Value = '''0.245'
$ws.Range('E31').Value = '  +4.79%  '
$ws.Range('D32').Value = '''0.177'
$ws.Range('E32').Value = '  +3.48%  '
$ws.Range('D33').Value = '''0.123'
$ws.Range('E33').Value = '  -9.63%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '''9.00'
$ws.Range('E34').Value = '  -2.36%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D35').Value = '''7.65'
$ws.Range('E35').Value = '  -5.68%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.154'
$ws.Range('E36').Value = '  -3.48%  '
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').Value = '''0.892'
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('D38').Value = '''25.31'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').Value = '''1.86'
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('D40').Value = '''23.99'
$ws.Range('D41').Value = '''0.437'
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''466.88'
$ws.Range('E42').Value = '  -5.39%  '
$ws.Range('B43').Value = 'MantraDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D43').Value = '''3.71'
$ws.Range('E43').Value = '  -3.95%  '
$ws.Range('D44').Value = '''1.25'
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  -8.85%  '
$ws.Range('D47').Value = '''160.88'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''0.669'
$ws.Range('E48').Value = '  -2.65%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '''1.82'
$ws.Range('E49').Value = '  -5.12%  '
$ws.Range('D50').Value = '''43.57'
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('D51').Value = '''0.998'
$ws.Range('E51').Value = '  -0.06%  '
